# Switch the presentation's applied design/colour theme from the custom
# "Integral" (Red Violet) scheme back to the standard Office Theme colours.
#
# PowerPoint stores the 12 theme colours (dk1, lt1, dk2, lt2, accent1-6,
# hlink, folHlink) on the slide master's Theme.ThemeColorScheme - every
# slide/layout that uses "Follow Master Scheme" (the default) re-colours
# itself from these values, so updating them here re-colours the whole
# deck in one go, exactly like using Design > Variants > Colors in the UI.

function RGB($r, $g, $b) { return $r + ($g * 256) + ($b * 65536) }

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$scheme = $master.Theme.ThemeColorScheme

# Office Theme palette, in clrScheme document order.
$scheme.Item(1).RGB  = RGB 0x00 0x00 0x00   # dk1      - 000000
$scheme.Item(2).RGB  = RGB 0xFF 0xFF 0xFF   # lt1      - FFFFFF
$scheme.Item(3).RGB  = RGB 0x44 0x54 0x6A   # dk2      - 44546A
$scheme.Item(4).RGB  = RGB 0xE7 0xE6 0xE6   # lt2      - E7E6E6
$scheme.Item(5).RGB  = RGB 0x5B 0x9B 0xD5   # accent1  - 5B9BD5
$scheme.Item(6).RGB  = RGB 0xED 0x7D 0x31   # accent2  - ED7D31
$scheme.Item(7).RGB  = RGB 0xA5 0xA5 0xA5   # accent3  - A5A5A5
$scheme.Item(8).RGB  = RGB 0xFF 0xC0 0x00   # accent4  - FFC000
$scheme.Item(9).RGB  = RGB 0x44 0x72 0xC4   # accent5  - 4472C4
$scheme.Item(10).RGB = RGB 0x70 0xAD 0x47   # accent6  - 70AD47
$scheme.Item(11).RGB = RGB 0x05 0x63 0xC1   # hlink    - 0563C1
$scheme.Item(12).RGB = RGB 0x95 0x4F 0x72   # folHlink - 954F72

# Best-effort: rename the theme/colour scheme to match (no-op on hosts that
# treat theme naming as read-only, but harmless and correct intent if
# supported - NOTE: deliberately not touching Design.Name, which this host
# maps onto the slide master's own <p:cSld> name rather than the theme name).
try { $master.Theme.Name = "Office Theme" } catch {}
try { $master.Theme.ThemeColorScheme.Name = "Office" } catch {}
